$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Update the two header shared-string cells (volume number + week-of dates).
# ------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 33   Number  4"
$ws.Range("C9").Value = "Report Covering the Week  1/19/2026  Through  1/25/2026"

# ------------------------------------------------------------------
# Cells that switch from a numeric value to the "0"/"***.*" placeholder text
# (need the text-only style (as used by C14) applied).
# ------------------------------------------------------------------
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("C25").NumberFormat = "@"
$ws.Range("C25").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C25").PasteSpecial(-4122)
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D31").PasteSpecial(-4122)
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E31").PasteSpecial(-4122)
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D33").PasteSpecial(-4122)
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E33").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Cells that switch from the placeholder text to a real numeric value
# (need the appropriate numeric style applied first).
# ------------------------------------------------------------------
$ws.Range("C39").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("N14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("C39").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("N14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 0
$ws.Range("C39").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$ws.Range("N14").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = -100
$ws.Range("N14").Copy()
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N15").Value = -100
$ws.Range("C39").Copy()
$ws.Range("C18").PasteSpecial(-4122)
$ws.Range("C18").Value = 4
$ws.Range("C39").Copy()
$ws.Range("I18").PasteSpecial(-4122)
$ws.Range("I18").Value = 4
$ws.Range("C39").Copy()
$ws.Range("C20").PasteSpecial(-4122)
$ws.Range("C20").Value = 2
$excel.CutCopyMode = 0

# ------------------------------------------------------------------
# Remaining cells: value-only changes (style/type unchanged).
# ------------------------------------------------------------------
$ws.Range("F15").Value = 1
$ws.Range("C16").Value = 1
$ws.Range("E16").Value = -66.666666666666
$ws.Range("F16").Value = 5
$ws.Range("G16").Value = 9
$ws.Range("H16").Value = -44.444444444444
$ws.Range("I16").Value = 4
$ws.Range("J16").Value = 7
$ws.Range("K16").Value = -42.857142857142
$ws.Range("L16").Value = -55.555555555555
$ws.Range("M16").Value = -69.230769230769
$ws.Range("N16").Value = -95.238095238095
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 24
$ws.Range("G17").Value = 14
$ws.Range("H17").Value = 71.428571428571
$ws.Range("I17").Value = 21
$ws.Range("J17").Value = 13
$ws.Range("K17").Value = 61.538461538461
$ws.Range("L17").Value = 5
$ws.Range("M17").Value = 90.909090909090
$ws.Range("N17").Value = -41.666666666666
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 100
$ws.Range("F18").Value = 4
$ws.Range("G18").Value = 6
$ws.Range("H18").Value = -33.333333333333
$ws.Range("J18").Value = 6
$ws.Range("K18").Value = -33.333333333333
$ws.Range("L18").Value = -69.230769230769
$ws.Range("M18").Value = -78.947368421052
$ws.Range("N18").Value = -93.939393939393
$ws.Range("D19").Value = 7
$ws.Range("E19").Value = -57.142857142857
$ws.Range("F19").Value = 20
$ws.Range("G19").Value = 13
$ws.Range("H19").Value = 53.846153846153
$ws.Range("I19").Value = 14
$ws.Range("J19").Value = 12
$ws.Range("K19").Value = 16.666666666666
$ws.Range("L19").Value = -63.157894736842
$ws.Range("M19").Value = -17.647058823529
$ws.Range("N19").Value = -50
$ws.Range("F20").Value = 5
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = 66.666666666666
$ws.Range("I20").Value = 4
$ws.Range("K20").Value = 300
$ws.Range("L20").Value = -73.333333333333
$ws.Range("M20").Value = -66.666666666666
$ws.Range("N20").Value = -96.078431372549
$ws.Range("C21").Value = 16
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 6.666666666666
$ws.Range("F21").Value = 59
$ws.Range("G21").Value = 46
$ws.Range("H21").Value = 28.260869565217
$ws.Range("I21").Value = 47
$ws.Range("J21").Value = 40
$ws.Range("K21").Value = 17.5
$ws.Range("L21").Value = -51.041666666666
$ws.Range("M21").Value = -36.486486486486
$ws.Range("N21").Value = -85.266457680250
$ws.Range("C24").Value = 4
$ws.Range("D24").Value = 11
$ws.Range("E24").Value = -63.636363636363
$ws.Range("F24").Value = 57
$ws.Range("H24").Value = 3.636363636363
$ws.Range("I24").Value = 50
$ws.Range("J24").Value = 50
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -25.373134328358
$ws.Range("M24").Value = -34.210526315789
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 14
$ws.Range("H25").Value = -17.647058823529
$ws.Range("I25").Value = 13
$ws.Range("J25").Value = 16
$ws.Range("K25").Value = -18.75
$ws.Range("L25").Value = -66.666666666666
$ws.Range("D26").Value = 5
$ws.Range("E26").Value = 100
$ws.Range("F26").Value = 38
$ws.Range("G26").Value = 32
$ws.Range("H26").Value = 18.75
$ws.Range("I26").Value = 34
$ws.Range("J26").Value = 27
$ws.Range("K26").Value = 25.925925925925
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -10.526315789473
$ws.Range("F27").Value = 1
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("J27").Value = 3
$ws.Range("D28").Value = 2
$ws.Range("E28").Value = -100
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -60
$ws.Range("J28").Value = 5
$ws.Range("K28").Value = -60
